$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 already has data in A7:F7 (MutationKey/Gene/case/sample/DNA+RNA
# sample) and T7:U7 (manual curation columns), but the per-variant-caller
# DNA columns (G:N) and RNA columns (O:S) were left blank for this
# variant. Fill them in with the "not called" placeholder text that is
# already used elsewhere in the sheet for the same situation.
$ws.Range("G7:N7").Value = "Not DNA Called"
$ws.Range("O7:S7").Value = "Not RNA Called"

# Reflect the newly entered range as the active selection, like Excel
# would after typing/filling this data in.
$ws.Range("G7:S7").Select() | Out-Null
